$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 844 (shifts existing rows 844:885 down to 845:886)
$ws.Rows.Item(844).Insert()

# Populate the new row 844 with the values from the diff.
# Column A holds a date written as plain text (e.g. "2026/02/21"), so force
# the cell to Text format first to avoid Excel auto-converting it to a date
# serial number.
$ws.Cells.Item(844, 1).NumberFormat = "@"
$ws.Cells.Item(844, 1).Value = "2026/02/21"
$ws.Cells.Item(844, 2).Value = "土"
$ws.Cells.Item(844, 3).Value = 20
$ws.Cells.Item(844, 4).Value = 201
